$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "RM 232" record entirely (originally row 26).
$ws.Rows.Item(26).Delete()

# Remove the "SC 92" record entirely.
# After the deletion above, the "SC 92" row has shifted up from 28 to 27.
$ws.Rows.Item(27).Delete()

# --- Update remaining cells to their new (re-imputed) values ---
$ws.Range("D2").Value = -13.5
$ws.Range("F4").Value = 17.97
$ws.Range("D12").Value = -14.1
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = 18.35
$ws.Range("D20").Value = -14
$ws.Range("D21").Value = -14.3
$ws.Range("F23").Value = 16.48
$ws.Range("F25").Value = 16.6
$ws.Range("C26").Value = 10.8
$ws.Range("C30").Value = 11.4
$ws.Range("D31").Value = -13.7
$ws.Range("D33").Value = -14.1

# --- Clear cells that have become missing values ---
$ws.Range("F3").ClearContents()
$ws.Range("F5").ClearContents()
$ws.Range("D6").ClearContents()
$ws.Range("F8").ClearContents()
$ws.Range("D14").ClearContents()
$ws.Range("F19").ClearContents()
$ws.Range("F22").ClearContents()
$ws.Range("D23").ClearContents()
$ws.Range("D24").ClearContents()
$ws.Range("C27").ClearContents()
$ws.Range("F27").ClearContents()
$ws.Range("C32").ClearContents()
